$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header cell G1, matching the style of the other header cells (row 1)
$ws.Range("G1").Value = "product_code"
$ws.Range("G1").Font.Bold = $true

# Apply the text number format to the whole product_code data column first
$ws.Range("G2:G14").NumberFormat = "@"

# Product code values, written in this specific order so the shared-string
# table is populated to match the source order.
$ws.Cells.Item(14, 7).Value = "5-1"
$ws.Cells.Item(8, 7).Value = "1-1"
$ws.Cells.Item(7, 7).Value = "1-3"
$ws.Cells.Item(5, 7).Value = "1-2"
$ws.Cells.Item(3, 7).Value = "1-4"
$ws.Cells.Item(4, 7).Value = "1-6"
$ws.Cells.Item(2, 7).Value = "1-9"
$ws.Cells.Item(9, 7).Value = "1-5"
$ws.Cells.Item(10, 7).Value = "1-7"
$ws.Cells.Item(11, 7).Value = "1-8"
$ws.Cells.Item(12, 7).Value = "1-10"
$ws.Cells.Item(13, 7).Value = "1-11"

$ws.Range("G14").Select()
